$wb = $excel.ActiveWorkbook

# Sheet "TestSuite" (first sheet): AdminPageTest RunMode changes from Y to N,
# selection moves from B3 to B6, and it is no longer the active/selected tab.
$wsSuite = $wb.Worksheets.Item("TestSuite")
$wsSuite.Range("B3").Value = "N"
$wsSuite.Range("B6").Select()

# Sheet "LoginPageTest" (second sheet): RunMode changes from Y to N,
# selection moves from C3 to D6, and it becomes the active/selected tab.
$wsLogin = $wb.Worksheets.Item("LoginPageTest")
$wsLogin.Range("C3").Value = "N"
$wsLogin.Activate()
$wsLogin.Range("D6").Select()
